# Generate Report for Handback
# Updates the localization-status report: mark zh-cn/de-de rows as handed
# back (in sync with en-US), refresh the "Latest Handback DateTime"
# timestamps, and clear the stale "Error Detail" message now that the
# handback files are up to date.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "Handed back: in sync with en-US"
$zhcn.Range("K2").Value = "2016-09-01 09:03:35"
$zhcn.Range("P2").Value = ""

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "Handed back: in sync with en-US"
$dede.Range("K2").Value = "2016-09-01 09:03:42"
$dede.Range("P2").Value = ""
